$wb = $excel.ActiveWorkbook

# The "Poland" sheet (last sheet) is the template for the new "UK" sheet.
$src = $wb.Worksheets.Item("Poland")

# Copy it to the end of the workbook.
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# Update the content for the UK market (order matches target shared-string order).
$newSheet.Range("B4").Value = "NGC-2741/T3365/T3367/T3364"
$newSheet.Range("B2").Value = "UK Market"

# Widen column B to fit the new content.
$newSheet.Columns.Item(2).ColumnWidth = 31

# Select B4 on the new sheet (becomes the active sheet/selection).
$newSheet.Range("B4").Select() | Out-Null
